$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Name column values
$ws.Range("A2").Value = "Goblin"
$ws.Range("A3").Value = "Grick"

# Remove the ATK/DEF columns (header + data) entirely
$ws.Range("B1:C3").Clear()

# Update selection to match target state
$ws.Range("B2:C3").Select()
